$p = $ppt.ActivePresentation

# Locate the paragraph that mentions the reminder-time example. It lives on
# the "Функционал" slide, in the body placeholder's bullet list (originally:
# Создание напоминаний с гибким парсингом времени (например, "каждый день в
# 8 утра").
$oldPhrase = ' времени (например, "каждый день в 8 утра").'
$newPhrase = ' времени (например, "каждый день в 8:00, через 5 минут, 1 июня в 15:00").'

$targetShape = $null
$targetSlide = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.HasText) {
                if ($shape.TextFrame.TextRange.Text.IndexOf($oldPhrase) -ge 0) {
                    $targetShape = $shape
                    $targetSlide = $slide
                }
            }
        }
    }
}

$tf = $targetShape.TextFrame
$tr = $tf.TextRange

# Find which paragraph holds the phrase.
$para1 = $null
for ($pi = 1; $pi -le $tr.Paragraphs().Count; $pi++) {
    $cand = $tr.Paragraphs($pi)
    if ($cand.Text.IndexOf($oldPhrase) -ge 0) {
        $para1 = $cand
    }
}

$fullParaText = $para1.Text
$startIdx0 = $fullParaText.IndexOf($oldPhrase)
$start = $startIdx0 + 1          # PowerPoint TextRange character positions are 1-based
$oldLen = $oldPhrase.Length

# Replace the trailing run's whole text first, keeping it a single run so the
# existing run formatting (Montserrat Medium / dk1 / 1800) carries over intact.
$whole = $para1.Characters($start, $oldLen)
$whole.Text = $newPhrase

# Now split that combined text back into four runs, matching the authored
# edit: "...в " | "8" | ":00, через 5 минут, 1 июня в 15:00" | "")."
$piece1 = ' времени (например, "каждый день в '
$piece2 = '8'
$piece3 = ':00, через 5 минут, 1 июня в 15:00'
$piece4 = '").'

$off = $start
$r1 = $para1.Characters($off, $piece1.Length); $off += $piece1.Length
$r2 = $para1.Characters($off, $piece2.Length); $off += $piece2.Length
$r3 = $para1.Characters($off, $piece3.Length); $off += $piece3.Length
$r4 = $para1.Characters($off, $piece4.Length); $off += $piece4.Length

# Re-assigning a font property to its own current value doesn't alter the
# visible formatting, but forces the engine to materialize each slice as an
# independent run (rather than silently re-merging it with its neighbour) -
# exactly what happened when the new runs were typed in the authoring app.
$r2.Font.Name = $r2.Font.Name
$r3.Font.Name = $r3.Font.Name
$r4.Font.Name = $r4.Font.Name
